$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Dataset sheet: drop the NamedThing/Person-style columns and replace with
# the new DCAT-style schema (was_generated_by, id, name, description).
# ---------------------------------------------------------------------------
$wsDataset = $wb.Worksheets.Item("Dataset")
$wsDataset.Range("D2:D1048576").Validation.Delete()

$wsDataset.Cells.Item(1,1).Value = "was_generated_by"
$wsDataset.Cells.Item(1,2).Value = "id"
$wsDataset.Cells.Item(1,3).Value = "name"
$wsDataset.Cells.Item(1,4).Value = "description"
$wsDataset.Cells.Item(1,5).Value = $null
$wsDataset.Cells.Item(1,6).Value = $null
$wsDataset.Cells.Item(1,7).Value = $null

# ---------------------------------------------------------------------------
# DatasetCollection -> renamed to Catalog, with a new column layout.
# ---------------------------------------------------------------------------
$wsCatalog = $wb.Worksheets.Item("DatasetCollection")
$wsCatalog.Name = "Catalog"

$wsCatalog.Cells.Item(1,1).Value = "has_dataset"
$wsCatalog.Cells.Item(1,2).Value = "was_generated_by"
$wsCatalog.Cells.Item(1,3).Value = "id"
$wsCatalog.Cells.Item(1,4).Value = "name"
$wsCatalog.Cells.Item(1,5).Value = "description"

# ---------------------------------------------------------------------------
# New sheet: Activity
# ---------------------------------------------------------------------------
$wsActivity = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsActivity.Name = "Activity"

$wsActivity.Cells.Item(1,1).Value = "type"
$wsActivity.Cells.Item(1,2).Value = "has_part"
$wsActivity.Cells.Item(1,3).Value = "had_object"
$wsActivity.Cells.Item(1,4).Value = "used_tool"

$wsActivity.Range("A2:A1048576").Validation.Add(3, 1, 1, '"data_curation,spectroscopy"')

# ---------------------------------------------------------------------------
# New sheet: ObjectOfInterest
# ---------------------------------------------------------------------------
$wsObjectOfInterest = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsObjectOfInterest.Name = "ObjectOfInterest"

$wsObjectOfInterest.Cells.Item(1,1).Value = "type"
$wsObjectOfInterest.Cells.Item(1,2).Value = "id"
$wsObjectOfInterest.Cells.Item(1,3).Value = "name"
$wsObjectOfInterest.Cells.Item(1,4).Value = "description"

$wsObjectOfInterest.Range("A2:A1048576").Validation.Add(3, 1, 1, '"cola"')

# ---------------------------------------------------------------------------
# New sheet: Tool
# ---------------------------------------------------------------------------
$wsTool = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsTool.Name = "Tool"

$wsTool.Cells.Item(1,1).Value = "type"
$wsTool.Cells.Item(1,2).Value = "has_part"
$wsTool.Cells.Item(1,3).Value = "id"
$wsTool.Cells.Item(1,4).Value = "name"
$wsTool.Cells.Item(1,5).Value = "description"

$wsTool.Range("A2:A1048576").Validation.Add(3, 1, 1, '"spectrometer"')
